# Add "Mobile_Number" (column E) to the Base Location sheet, matching the
# final-submission dashboard export: new header cell, 103 phone numbers,
# and the font/number-format touches that came along with the pasted-in data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: E1 gets the same bold header style as the other headers ---
$ws.Range("D1").Copy($ws.Range("E1")) | Out-Null
$ws.Range("E1").Value = "Mobile_Number"

# --- Custom cell style carried over from the source data (Arial 10) ---
$customStyle = $wb.Styles.Add("Normal 2 2")
$customStyle.Font.Name = "Arial"
$customStyle.Font.Size = 10

# --- Column E data: row|styleGroup|value ---
# styleGroup: 0 = default, 2/4 = Arial 10 ("Normal 2 2"), 3 = Calibri + "0" number format
$data = @(
    "2|0|9693888954",
    "3|0|9835214508",
    "4|0|9334335339",
    "5|0|8971709801",
    "6|0|8294447787",
    "7|0|8966324512",
    "8|0|9894715876",
    "9|0|8971709213",
    "10|0|9693888741",
    "11|0|8971702564",
    "12|0|8845369871",
    "13|0|9835456012",
    "14|2|9934164003",
    "15|2|9934164004",
    "16|2|9934164005",
    "17|2|9934164027",
    "18|2|9934164032",
    "19|0|9216000081",
    "20|0|9216001808",
    "21|0|9216001859",
    "22|0|9216103705",
    "23|0|9755974116",
    "24|0|9755858671",
    "25|0|9755759990",
    "26|0|9755369104",
    "27|0|9993558512",
    "28|0|9993558523",
    "29|0|9993559112",
    "30|0|9993558365",
    "31|0|9993559074",
    "32|3|9899570762",
    "33|3|9899570764",
    "34|3|9899570808",
    "35|3|9899570813",
    "36|0|9825269655",
    "37|0|9824060573",
    "38|0|9825865168",
    "39|0|9825022295",
    "40|0|9824050801",
    "41|2|9898300015",
    "42|2|9898300028",
    "43|2|9898300040",
    "44|2|9898300041",
    "45|0|9210659890",
    "46|0|9213244346",
    "47|0|9213244678",
    "48|0|9213359351",
    "49|0|9213222830",
    "50|0|9213820352",
    "51|0|9213352603",
    "52|0|9210306856",
    "53|0|9210458065",
    "54|0|9210452310",
    "55|0|9210580412",
    "56|0|9213601721",
    "57|2|9827700015",
    "58|2|9827700092",
    "59|2|9827708219",
    "60|2|9827708258",
    "61|2|9827708269",
    "62|4|9444307544",
    "63|4|9841184725",
    "64|4|9444032428",
    "65|4|9841518221",
    "66|4|9840231471",
    "67|3|9899570808",
    "68|3|9899570813",
    "69|0|9825269655",
    "70|0|9824060573",
    "71|0|9825865168",
    "72|0|9825022295",
    "73|0|9824050801",
    "74|2|9898300015",
    "75|2|9898300028",
    "76|2|9898300040",
    "77|2|9898300041",
    "78|0|9210659890",
    "79|0|9213244346",
    "80|0|9213244678",
    "81|0|9213359351",
    "82|0|8971709801",
    "83|0|8294447787",
    "84|0|8966324512",
    "85|0|9894715876",
    "86|0|8971709213",
    "87|0|9693888741",
    "88|0|8971702564",
    "89|0|8845369871",
    "90|0|9835456012",
    "91|2|9934164003",
    "92|2|9934164004",
    "93|2|9934164005",
    "94|2|9934164027",
    "95|2|9934164032",
    "96|0|9216001808",
    "97|0|9216001859",
    "98|0|9216103705",
    "99|0|9755974116",
    "100|0|9755858671",
    "101|0|9755759990",
    "102|0|9755369104",
    "103|0|9993558512",
    "104|0|9993558523"
)

foreach ($entry in $data) {
    $parts = $entry.Split("|")
    $r = [int]$parts[0]
    $styleGroup = $parts[1]
    $val = [double]$parts[2]
    $cell = $ws.Cells.Item($r, 5)
    $cell.Value = $val
    if ($styleGroup -eq "2" -or $styleGroup -eq "4") {
        $cell.Style = "Normal 2 2"
    } elseif ($styleGroup -eq "3") {
        $cell.NumberFormat = "0"
    }
}

# --- Column E width (matches the pasted data's source column width) ---
$ws.Columns.Item(5).ColumnWidth = 14.14

# --- View: selection moves to E3, scrolled back to top-left ---
$ws.Range("E3").Select()
